# Update the per-tab SQL query cells (column B) on Sheet1 to reflect the
# corrected/updated queries (column references renamed, extra WHERE
# clauses/CASE branches added), then leave the selection on the last
# edited cell (B5), matching the author's final saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ProgramsTab query: "Focus Area" -> "Special Topic" label, and the
# Data Location Details CASE now falls back to prg.program_acronym instead
# of prg.website.
$ws.Range("B2").Value = 'SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Special Topic",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.program_acronym     
        ELSE prg.data_link
    END AS "Data Location Details" 
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE ''%Multiple Myeloma%''
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;'

# Row 3 - ProjectsTab query: prj.org_name -> prj.project_org_name
$ws.Range("B3").Value = 'SELECT DISTINCT
    prj.project_id AS "Project ID", 
    prj.project_title AS "Project Title",
    prj.project_org_name AS "Organization",
    prj.project_start_date AS "Project Start Date",
    prj.project_end_date AS "Project End Date"
FROM 
    df_project prj
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
     prg.cancer_type LIKE ''%Multiple Myeloma%''
ORDER BY 
    lower(prj.project_id) ASC
LIMIT 100;'

# Row 4 - GrantsTab query: gnt.project_end_date -> gnt.grant_end_date
$ws.Range("B4").Value = 'SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.grant_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type  LIKE ''%Multiple Myeloma%''
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;'

# Row 5 - PublicationsTab query: pub.title -> pub.publication_title, plus two
# additional CASE WHEN branches for relative_citation_ratio 1.0 and 5.0.
$ws.Range("B5").Value = 'SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.publication_title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN ''0''
    WHEN pub.relative_citation_ratio = 7.0 THEN ''7''
    WHEN pub.relative_citation_ratio = 2.0 THEN ''2''
  WHEN pub.relative_citation_ratio = 1.0 THEN ''1''
WHEN pub.relative_citation_ratio = 5.0 THEN ''5''
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
     prg.cancer_type  LIKE ''%Multiple Myeloma%''
ORDER BY 
    lower(pub.pmid) ASC
LIMIT 100;'

# Match the author's final cursor position/selection in the saved file.
$ws.Range("B5").Select()
